$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  A=34957; B="Lucas Moreira";            C="Juridico";          D="Consulta medica";   E=6; F=45097; G=3750.81}
    @{Row=3;  A=65256; B="Alícia Pacheco";            C="Financeiro";        D="Doenca";             E=4; F=45089; G=5854.66}
    @{Row=4;  A=58745; B="Dr. Henrique Montenegro";   C="Marketing";         D="Doenca";             E=5; F=45106; G=2041.96}
    @{Row=5;  A=6782;  B="Samuel Gomes";              C="Recursos Humanos";  D="Problemas pessoais"; E=2; F=45101; G=9437.139999999999}
    @{Row=6;  A=39041; B="Julia Martins";             C="P&D";               D="Outros";             E=4; F=45094; G=6212.28}
    @{Row=7;  A=38566; B="Gael Novaes";               C="Recursos Humanos";  D="Problemas pessoais"; E=6; F=45085; G=7976.82}
    @{Row=8;  A=87888; B="Yuri da Costa";             C="Financeiro";        D="Consulta medica";   E=5; F=45104; G=6314.42}
    @{Row=9;  A=91613; B="Enzo Moura";                C="Recursos Humanos";  D="Consulta medica";   E=7; F=45083; G=4044.01}
    @{Row=10; A=15378; B="Ravi Lucca Cassiano";       C="Engenharia";        D="Doenca";             E=2; F=45086; G=9880.23}
    @{Row=11; A=87327; B="Antônio da Rocha";          C="Recursos Humanos";  D="Outros";             E=1; F=45078; G=2038.75}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}
